# "update: excel date update"
#
# 1) Insert a new "Description" sheet in front of "Funding" with a textbox
#    that explains the funding-vs-expense timing mismatch.
# 2) Convert the text-formatted date columns on "Funding" (Valid To) and
#    "Expense" (Latest Payment Date) into real date values, and move every
#    date column onto the custom d/mm/yyyy;@ number format.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New "Description" sheet, inserted before "Funding", with a summary textbox.
# Do this first so the later Worksheets.Item(...) lookups are never stale.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("Funding")
$description = $wb.Worksheets.Add($beforeSheet)
$description.Name = "Description"

$q = [char]0x2019
$rupee = [char]0x20B9
$nl = [char]10

$text = "What" + $q + "s happening: Most funding expires early, but most expenses happen later. Funds expiring by June " + $rupee + "60k" + $nl + "Expenses in Aug" + [char]0x2013 + "Nov " + $rupee + "70k" + $nl + "Late-year funds only " + $rupee + "6k" + $nl + " Early expenses look fine, later ones pile up uncovered."

$box = $description.Shapes.AddTextbox(1, 13.2, 12.6, 192.0, 73.0)
$box.Name = "TextBox 1"
$box.TextFrame.Characters().Text = $text

# ---------------------------------------------------------------------------
# Fetch fresh handles by name now that the sheet collection has shifted.
# ---------------------------------------------------------------------------
$funding = $wb.Worksheets.Item("Funding")
$expense = $wb.Worksheets.Item("Expense")

# Funding sheet: "Valid From" (D) keeps its values, just gets the new custom
# date format. "Valid To" (E) was stored as literal text ("6/30/2025" /
# "12/31/2025") - replace with real date serials and the new format.
$funding.Range("D2:D16").NumberFormat = "d/mm/yyyy;@"

$funding.Range("E2:E11").Value = 45838
$funding.Range("E12:E16").Value = 46022
$funding.Range("E2:E16").NumberFormat = "d/mm/yyyy;@"

# Expense sheet: "Latest Payment Date" (F) is a mix of real date serials and
# literal text dates - normalize everything to real dates + new format.
$expense.Range("F4").Value = 45741
$expense.Range("F7").Value = 45894
$expense.Range("F10").Value = 45955
$expense.Range("F11").Value = 45981
$expense.Range("F2:F11").NumberFormat = "d/mm/yyyy;@"

# ---------------------------------------------------------------------------
# Restore per-sheet selections, then land back on the Description sheet.
# ---------------------------------------------------------------------------
$funding.Activate()
$funding.Range("E24").Select()
$expense.Activate()
$expense.Range("I8").Select()
$description.Activate()
